$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, matching the style of the other header cells (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Fill the new "Save" column (H) for each data row with 0, except row 3 which is 1
for ($r = 2; $r -le 69; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}

$ws.Cells.Item(3, 8).Value = 1
